$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.237.58'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '2.621.24'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.05'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.39'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.552'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.46%  '
$ws.Range('D9').Value = '2.619.99'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.347'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').Value = '3.102.44'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').Value = '67.264.13'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').Value = '2.624.77'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.10'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '362.72'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.48'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.55%  '
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('E23').Value = '  +3.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.92'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.05'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('D27').Value = '2.762.05'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '578.34'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.60%  '
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.79'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('E33').Value = '  -0.98%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.127'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.51'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.87'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '156.90'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.11'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.24'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.80'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.14'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '156.24'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('E48').Value = '  -2.88%  '
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.621'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.53'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.32%  '
